$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1428654.9
$ws.Range("I9").Value = 1666738.6
$ws.Range("J9").Value = 152
$ws.Range("K9").Value = 1666738.6
$ws.Range("L9").Value = 152
$ws.Range("M9").Value = -1666569.6
$ws.Range("N9").Value = -490
$ws.Range("H18").Value = 927.1429000000001
$ws.Range("I18").Value = 498.33334
$ws.Range("K18").Value = 498.33334
$ws.Range("M18").Value = -214.33334
$ws.Range("H21").Value = 31111
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 31111
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 31111
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -32047
$ws.Range("H23").Value = 31111
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 31111
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 31111
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -31579
$ws.Range("H29").Value = 166669660
$ws.Range("I29").Value = 500000000
$ws.Range("J29").Value = 4501.5
$ws.Range("K29").Value = 1500000000
$ws.Range("L29").Value = 13504.5
$ws.Range("M29").Value = -1499999719
$ws.Range("N29").Value = -14066.5
$ws.Range("H32").Value = 2498.5454
$ws.Range("I32").Value = 1900
$ws.Range("J32").Value = 2631.5557
$ws.Range("K32").Value = 1900
$ws.Range("L32").Value = 2631.5557
$ws.Range("M32").Value = -1574
$ws.Range("N32").Value = -3283.5557
$ws.Range("H34").Value = 1259.7778
$ws.Range("I34").Value = 1259.7778
$ws.Range("K34").Value = 1259.7778
$ws.Range("M34").Value = -1056.7778
$ws.Range("H36").Value = 1259.7778
$ws.Range("I36").Value = 1259.7778
$ws.Range("K36").Value = 1259.7778
$ws.Range("M36").Value = -544.7778000000001
$ws.Range("H40").Value = 1721.3055
$ws.Range("I40").Value = 1478.6
$ws.Range("J40").Value = 2272.9092
$ws.Range("K40").Value = 1478.6
$ws.Range("L40").Value = 2272.9092
$ws.Range("M40").Value = -1303.6
$ws.Range("N40").Value = -2622.9092
$ws.Range("H43").Value = 9045
$ws.Range("I43").Value = 7569.857
$ws.Range("J43").Value = 10077.6
$ws.Range("K43").Value = 7569.857
$ws.Range("L43").Value = 10077.6
$ws.Range("M43").Value = -7500.857
$ws.Range("N43").Value = -10215.6
$ws.Range("H62").Value = 38981.582
$ws.Range("I62").Value = 46082
$ws.Range("K62").Value = 46082
$ws.Range("M62").Value = -45458
$ws.Range("H64").Value = 7758.3613
$ws.Range("I64").Value = 3778
$ws.Range("J64").Value = 9085.147999999999
$ws.Range("K64").Value = 3778
$ws.Range("L64").Value = 9085.147999999999
$ws.Range("M64").Value = -3530
$ws.Range("N64").Value = -9581.147999999999
$ws.Range("H65").Value = 38981.582
$ws.Range("I65").Value = 46082
$ws.Range("K65").Value = 230410
$ws.Range("M65").Value = -227290
$ws.Range("H67").Value = 7758.3613
$ws.Range("I67").Value = 3778
$ws.Range("J67").Value = 9085.147999999999
$ws.Range("K67").Value = 3778
$ws.Range("L67").Value = 9085.147999999999
$ws.Range("M67").Value = -2920
$ws.Range("N67").Value = -10801.148
$ws.Range("H86").Value = 2553.3215
$ws.Range("I86").Value = 3174.389
$ws.Range("J86").Value = 1435.4
$ws.Range("K86").Value = 3174.389
$ws.Range("L86").Value = 1435.4
$ws.Range("M86").Value = -2051.389
$ws.Range("N86").Value = -3681.4
$ws.Range("H89").Value = 2553.3215
$ws.Range("I89").Value = 3174.389
$ws.Range("J89").Value = 1435.4
$ws.Range("K89").Value = 15871.945
$ws.Range("L89").Value = 7177
$ws.Range("M89").Value = -10255.945
$ws.Range("N89").Value = -18409
$ws.Range("H106").Value = 5090.6113
$ws.Range("I106").Value = 4545.0713
$ws.Range("J106").Value = 7000
$ws.Range("K106").Value = 4545.0713
$ws.Range("L106").Value = 7000
$ws.Range("M106").Value = -3914.0713
$ws.Range("N106").Value = -8262
$ws.Range("H107").Value = 913.61536
$ws.Range("I107").Value = 913.61536
$ws.Range("K107").Value = 913.61536
$ws.Range("M107").Value = 1006.38464
$ws.Range("H111").Value = 3391.0557
$ws.Range("I111").Value = 3434
$ws.Range("K111").Value = 10302
$ws.Range("M111").Value = -7235
$ws.Range("H116").Value = 2948.25
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H127").Value = 1358.1818
$ws.Range("I127").Value = 1339
$ws.Range("J127").Value = 1444.5
$ws.Range("K127").Value = 4017
$ws.Range("L127").Value = 4333.5
$ws.Range("M127").Value = 943
$ws.Range("N127").Value = -14253.5
$ws.Range("H129").Value = 1862.4445
$ws.Range("I129").Value = 1303.3334
$ws.Range("J129").Value = 2980.6667
$ws.Range("K129").Value = 3910.0002
$ws.Range("L129").Value = 8942.000100000001
$ws.Range("M129").Value = 1089.9998
$ws.Range("N129").Value = -18942.0001
$ws.Range("H131").Value = 9935
$ws.Range("I131").Value = 8933.666999999999
$ws.Range("J131").Value = 10310.5
$ws.Range("K131").Value = 26801.001
$ws.Range("L131").Value = 30931.5
$ws.Range("M131").Value = -21761.001
$ws.Range("N131").Value = -41011.5
$ws.Range("H132").Value = 1789.8235
$ws.Range("I132").Value = 1223.25
$ws.Range("J132").Value = 3149.6
$ws.Range("K132").Value = 3669.75
$ws.Range("L132").Value = 9448.799999999999
$ws.Range("M132").Value = -1139.75
$ws.Range("N132").Value = -14508.8
$ws.Range("H135").Value = 678.1786
$ws.Range("I135").Value = 692.1852
$ws.Range("K135").Value = 6229.6668
$ws.Range("M135").Value = -3694.6668
$ws.Range("H137").Value = 1226.6571
$ws.Range("I137").Value = 1115.6072
$ws.Range("J137").Value = 1670.8572
$ws.Range("K137").Value = 3346.8216
$ws.Range("L137").Value = 5012.571599999999
$ws.Range("M137").Value = -796.8215999999998
$ws.Range("N137").Value = -10112.5716
$ws.Range("H138").Value = 2832.6553
$ws.Range("I138").Value = 1746.1111
$ws.Range("K138").Value = 5238.3333
$ws.Range("M138").Value = -98.33330000000024
$ws.Range("H141").Value = 6635.5713
$ws.Range("I141").Value = 7737.5
$ws.Range("J141").Value = 5166.3335
$ws.Range("K141").Value = 23212.5
$ws.Range("L141").Value = 15499.0005
$ws.Range("M141").Value = -18032.5
$ws.Range("N141").Value = -25859.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3753.1353
$ws.Range("I32").Value = 3709.4333
$ws.Range("J32").Value = 3940.4285
$ws.Range("K32").Value = 3709.4333
$ws.Range("L32").Value = 3940.4285
$ws.Range("M32").Value = -3422.4333
$ws.Range("N32").Value = -4514.4285
$ws.Range("H37").Value = 34728
$ws.Range("J37").Value = 34728
$ws.Range("L37").Value = 34728
$ws.Range("N37").Value = -35274
$ws.Range("H53").Value = 10039
$ws.Range("I53").Value = 10039
$ws.Range("K53").Value = 10039
$ws.Range("M53").Value = -9357
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9258
$ws.Range("H61").Value = 4724
$ws.Range("I61").Value = 4366.84
$ws.Range("J61").Value = 5999.5713
$ws.Range("K61").Value = 4366.84
$ws.Range("L61").Value = 5999.5713
$ws.Range("M61").Value = -4154.84
$ws.Range("N61").Value = -6423.5713
$ws.Range("H63").Value = 5267.5713
$ws.Range("I63").Value = 4772.4287
$ws.Range("K63").Value = 4772.4287
$ws.Range("M63").Value = -4086.4287
$ws.Range("H66").Value = 5267.5713
$ws.Range("I66").Value = 4772.4287
$ws.Range("K66").Value = 23862.1435
$ws.Range("M66").Value = -20430.1435
$ws.Range("H74").Value = 2134.8823
$ws.Range("I74").Value = 1819.5333
$ws.Range("K74").Value = 1819.5333
$ws.Range("M74").Value = -945.5333000000001
$ws.Range("H77").Value = 2134.8823
$ws.Range("I77").Value = 1819.5333
$ws.Range("K77").Value = 9097.666499999999
$ws.Range("M77").Value = -4729.666499999999
$ws.Range("H109").Value = 65000
$ws.Range("J109").Value = 65000
$ws.Range("L109").Value = 65000
$ws.Range("N109").Value = -67774
$ws.Range("H122").Value = 2657.5557
$ws.Range("I122").Value = 2573
$ws.Range("K122").Value = 7719
$ws.Range("M122").Value = -5269
$ws.Range("H132").Value = 1490
$ws.Range("I132").Value = 1485.1428
$ws.Range("K132").Value = 4455.428400000001
$ws.Range("M132").Value = -1925.428400000001
$ws.Range("H136").Value = 4724
$ws.Range("I136").Value = 4366.84
$ws.Range("J136").Value = 5999.5713
$ws.Range("K136").Value = 13100.52
$ws.Range("L136").Value = 17998.7139
$ws.Range("M136").Value = -10550.52
$ws.Range("N136").Value = -23098.7139
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1619.5834
$ws.Range("I20").Value = 1167.4
$ws.Range("J20").Value = 1942.5714
$ws.Range("K20").Value = 1167.4
$ws.Range("L20").Value = 1942.5714
$ws.Range("M20").Value = -920.4000000000001
$ws.Range("N20").Value = -2436.5714
$ws.Range("H86").Value = 1268416.6
$ws.Range("I86").Value = 1839.1111
$ws.Range("J86").Value = 3168282.8
$ws.Range("K86").Value = 1839.1111
$ws.Range("L86").Value = 3168282.8
$ws.Range("M86").Value = -716.1111000000001
$ws.Range("N86").Value = -3170528.8
$ws.Range("H89").Value = 1268416.6
$ws.Range("I89").Value = 1839.1111
$ws.Range("J89").Value = 3168282.8
$ws.Range("K89").Value = 9195.5555
$ws.Range("L89").Value = 15841414
$ws.Range("M89").Value = -3579.5555
$ws.Range("N89").Value = -15852646
$ws.Range("H94").Value = 2200.3
$ws.Range("I94").Value = 1417.1666
$ws.Range("K94").Value = 1417.1666
$ws.Range("M94").Value = -966.1666
$ws.Range("H99").Value = 3221.7368
$ws.Range("I99").Value = 1503
$ws.Range("J99").Value = 4015
$ws.Range("K99").Value = 1503
$ws.Range("L99").Value = 4015
$ws.Range("M99").Value = -5
$ws.Range("N99").Value = -7011
$ws.Range("H105").Value = 2014.5
$ws.Range("I105").Value = 1890.95
$ws.Range("K105").Value = 1890.95
$ws.Range("M105").Value = -143.95
$ws.Range("H134").Value = 1750.1395
$ws.Range("I134").Value = 1694.0731
$ws.Range("J134").Value = 2899.5
$ws.Range("K134").Value = 5082.219300000001
$ws.Range("L134").Value = 8698.5
$ws.Range("M134").Value = -2547.219300000001
$ws.Range("N134").Value = -13768.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 26051.75
$ws.Range("I41").Value = 5875
$ws.Range("J41").Value = 46228.5
$ws.Range("K41").Value = 5875
$ws.Range("L41").Value = 46228.5
$ws.Range("M41").Value = -5447
$ws.Range("N41").Value = -47084.5
$ws.Range("H58").Value = 2715.7334
$ws.Range("I58").Value = 3112.4443
$ws.Range("J58").Value = 2120.6667
$ws.Range("K58").Value = 3112.4443
$ws.Range("L58").Value = 2120.6667
$ws.Range("M58").Value = -2909.4443
$ws.Range("N58").Value = -2526.6667
$ws.Range("H99").Value = 3778.125
$ws.Range("I99").Value = 3627.75
$ws.Range("J99").Value = 3928.5
$ws.Range("K99").Value = 3627.75
$ws.Range("L99").Value = 3928.5
$ws.Range("M99").Value = -2129.75
$ws.Range("N99").Value = -6924.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3778.125
$ws.Range("I126").Value = 3627.75
$ws.Range("J126").Value = 3928.5
$ws.Range("K126").Value = 10883.25
$ws.Range("L126").Value = 11785.5
$ws.Range("M126").Value = -8413.25
$ws.Range("N126").Value = -16725.5
$ws.Range("H132").Value = 1567.5
$ws.Range("I132").Value = 1487.381
$ws.Range("J132").Value = 1904
$ws.Range("K132").Value = 4462.143
$ws.Range("L132").Value = 5712
$ws.Range("M132").Value = -1932.143
$ws.Range("N132").Value = -10772
$ws.Range("H134").Value = 1804.7727
$ws.Range("I134").Value = 1767.2778
$ws.Range("J134").Value = 1973.5
$ws.Range("K134").Value = 5301.8334
$ws.Range("L134").Value = 5920.5
$ws.Range("M134").Value = -2766.8334
$ws.Range("N134").Value = -10990.5
$ws.Range("H136").Value = 2715.7334
$ws.Range("I136").Value = 3112.4443
$ws.Range("J136").Value = 2120.6667
$ws.Range("K136").Value = 9337.332900000001
$ws.Range("L136").Value = 6362.000100000001
$ws.Range("M136").Value = -6787.332900000001
$ws.Range("N136").Value = -11462.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 21.8
$ws.Range("I40").Value = 21.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 87.2
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -18.2
$ws.Range("N40").ClearContents()
$ws.Range("H58").Value = 500
$ws.Range("I58").Value = 500
$ws.Range("K58").Value = 1500
$ws.Range("M58").Value = -1372
$ws.Range("H80").Value = 1266.6666
$ws.Range("I80").Value = 2800
$ws.Range("K80").Value = 8400
$ws.Range("M80").Value = -7464
$ws.Range("H83").Value = 1266.6666
$ws.Range("I83").Value = 2800
$ws.Range("K83").Value = 25200
$ws.Range("M83").Value = -20520
$ws.Range("H88").Value = 24555.555
$ws.Range("J88").Value = 26384.615
$ws.Range("L88").Value = 79153.845
$ws.Range("N88").Value = -80009.845
$ws.Range("H91").Value = 24555.555
$ws.Range("J91").Value = 26384.615
$ws.Range("L91").Value = 79153.845
$ws.Range("N91").Value = -82117.845
$ws.Range("H99").Value = 15598.917
$ws.Range("I99").Value = 7148.625
$ws.Range("K99").Value = 21445.875
$ws.Range("M99").Value = -19199.875
$ws.Range("H103").Value = 587.0909
$ws.Range("I103").Value = 172.14285
$ws.Range("J103").Value = 1313.25
$ws.Range("K103").Value = 516.4285500000001
$ws.Range("L103").Value = 3939.75
$ws.Range("M103").Value = 362.5714499999999
$ws.Range("N103").Value = -5697.75
$ws.Range("H132").Value = 1663.0952
$ws.Range("I132").Value = 993.2727
$ws.Range("J132").Value = 2399.9
$ws.Range("K132").Value = 8939.454299999999
$ws.Range("L132").Value = 21599.1
$ws.Range("M132").Value = -6409.454299999999
$ws.Range("N132").Value = -26659.1

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 63.25
$ws.Range("J2").Value = 157.8
$ws.Range("L2").Value = 157.8
$ws.Range("N2").Value = -383.8
$ws.Range("H11").Value = 11869583
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10278
$ws.Range("H70").Value = 22399.8
$ws.Range("I70").Value = 22399.8
$ws.Range("K70").Value = 22399.8
$ws.Range("M70").Value = -22129.8
$ws.Range("H73").Value = 22399.8
$ws.Range("I73").Value = 22399.8
$ws.Range("K73").Value = 22399.8
$ws.Range("M73").Value = -21463.8
$ws.Range("H80").Value = 6034.1816
$ws.Range("I80").Value = 4396
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 4396
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = -3398
$ws.Range("N80").Value = -9996
$ws.Range("H83").Value = 6034.1816
$ws.Range("I83").Value = 4396
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 21980
$ws.Range("L83").Value = 40000
$ws.Range("M83").Value = -16988
$ws.Range("N83").Value = -49984
$ws.Range("H102").Value = 3564.204
$ws.Range("I102").Value = 2385.5151
$ws.Range("J102").Value = 5995.25
$ws.Range("K102").Value = 2385.5151
$ws.Range("L102").Value = 5995.25
$ws.Range("M102").Value = -763.5151000000001
$ws.Range("N102").Value = -9239.25
$ws.Range("H113").Value = 8868.235000000001
$ws.Range("I113").Value = 6822.857
$ws.Range("J113").Value = 10300
$ws.Range("K113").Value = 6822.857
$ws.Range("L113").Value = 10300
$ws.Range("M113").Value = -4652.857
$ws.Range("N113").Value = -14640
$ws.Range("H126").Value = 3333.3333
$ws.Range("I126").Value = 3333.3333
$ws.Range("K126").Value = 9999.999899999999
$ws.Range("M126").Value = -7529.999899999999
$ws.Range("H132").Value = 1843.84
$ws.Range("I132").Value = 1556.0952
$ws.Range("K132").Value = 4668.2856
$ws.Range("M132").Value = -2138.2856
$ws.Range("H139").Value = 99473.25
$ws.Range("J139").Value = 99473.25
$ws.Range("L139").Value = 99473.25
$ws.Range("N139").Value = -109753.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 11875
$ws.Range("J18").Value = 11875
$ws.Range("L18").Value = 11875
$ws.Range("N18").Value = -12219
$ws.Range("H20").Value = 933.3333
$ws.Range("I20").Value = 900
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -674
$ws.Range("N20").Value = -1452
$ws.Range("H22").Value = 3521.64
$ws.Range("I22").Value = 3117
$ws.Range("K22").Value = 3117
$ws.Range("M22").Value = -2822
$ws.Range("H27").Value = 3521.64
$ws.Range("I27").Value = 3117
$ws.Range("K27").Value = 3117
$ws.Range("M27").Value = -3010
$ws.Range("H32").Value = 20000
$ws.Range("I32").Value = 20000
$ws.Range("K32").Value = 20000
$ws.Range("M32").Value = -19683
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H55").Value = 1885.25
$ws.Range("I55").Value = 310.14285
$ws.Range("J55").Value = 3110.3333
$ws.Range("K55").Value = 310.14285
$ws.Range("L55").Value = 3110.3333
$ws.Range("M55").Value = -137.14285
$ws.Range("N55").Value = -3456.3333
$ws.Range("H61").Value = 3791.6538
$ws.Range("I61").Value = 2661.8125
$ws.Range("K61").Value = 2661.8125
$ws.Range("M61").Value = -2459.8125
$ws.Range("H82").Value = 2107.2173
$ws.Range("I82").Value = 872.75
$ws.Range("J82").Value = 4928.857
$ws.Range("K82").Value = 872.75
$ws.Range("L82").Value = 4928.857
$ws.Range("M82").Value = -511.75
$ws.Range("N82").Value = -5650.857
$ws.Range("H85").Value = 2107.2173
$ws.Range("I85").Value = 872.75
$ws.Range("J85").Value = 4928.857
$ws.Range("K85").Value = 872.75
$ws.Range("L85").Value = 4928.857
$ws.Range("M85").Value = 375.25
$ws.Range("N85").Value = -7424.857
$ws.Range("H113").Value = 3791.6538
$ws.Range("I113").Value = 2661.8125
$ws.Range("K113").Value = 2661.8125
$ws.Range("M113").Value = -491.8125
$ws.Range("H122").Value = 7638.467
$ws.Range("I122").Value = 7505.5
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 22516.5
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -20066.5
$ws.Range("N122").Value = -33400
$ws.Range("H132").Value = 2934.9375
$ws.Range("I132").Value = 2611.3572
$ws.Range("K132").Value = 7834.071599999999
$ws.Range("M132").Value = -5304.071599999999
$ws.Range("H136").Value = 20755.072
$ws.Range("I136").Value = 1932.3334
$ws.Range("J136").Value = 34872.125
$ws.Range("K136").Value = 5797.0002
$ws.Range("L136").Value = 104616.375
$ws.Range("M136").Value = -3247.0002
$ws.Range("N136").Value = -109716.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1266.6666
$ws.Range("I113").Value = 984.7692
$ws.Range("J113").Value = 1724.75
$ws.Range("K113").Value = 2954.3076
$ws.Range("L113").Value = 5174.25
$ws.Range("M113").Value = -784.3076000000001
$ws.Range("N113").Value = -9514.25
$ws.Range("H122").Value = 7894.25
$ws.Range("I122").Value = 3577.5
$ws.Range("J122").Value = 9333.166999999999
$ws.Range("K122").Value = 10732.5
$ws.Range("L122").Value = 27999.501
$ws.Range("M122").Value = -8282.5
$ws.Range("N122").Value = -32899.501
$ws.Range("H123").Value = 44749.5
$ws.Range("J123").Value = 44749.5
$ws.Range("L123").Value = 44749.5
$ws.Range("N123").Value = -54549.5
$ws.Range("H126").Value = 2668.7144
$ws.Range("I126").Value = 2545.5386
$ws.Range("J126").Value = 2868.875
$ws.Range("K126").Value = 7636.6158
$ws.Range("L126").Value = 8606.625
$ws.Range("M126").Value = -5166.6158
$ws.Range("N126").Value = -13546.625
$ws.Range("H132").Value = 2297.02
$ws.Range("I132").Value = 2090.7144
$ws.Range("K132").Value = 6272.1432
$ws.Range("M132").Value = -3742.1432
$ws.Range("J139").Value = 130000
$ws.Range("L139").Value = 130000
$ws.Range("N139").Value = -140280
